$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet from "Sheet1" to "Login"
$ws.Name = "Login"

# Update the header row text (row 1)
$ws.Range("A1").Value = "Username"
$ws.Range("B1").Value = "Password"
$ws.Range("C1").Value = "ExpectedResult"

# Update the data rows (rows 2-3)
$ws.Range("A2").Value = "standard_user"
$ws.Range("B2").Value = "secret_sauce"
$ws.Range("C2").Value = "Pass"

$ws.Range("A3").Value = "locked_out_user"
$ws.Range("B3").Value = "secret_sauce"
$ws.Range("C3").Value = "Fail"

# Header formatting: Arial, bold (already bold), light-grey fill, medium light-grey border
$header = $ws.Range("A1:C1")
$header.Font.Name = "Arial"
$header.Interior.Color = 15921906
$header.Borders.Weight = -4138
$header.Borders.Color = 14540253

# Body formatting: Arial, medium light-grey border
$body = $ws.Range("A2:C3")
$body.Font.Name = "Arial"
$body.Borders.Weight = -4138
$body.Borders.Color = 14540253

# Update the active selection shown in the sheet view
$null = $ws.Range("B13").Select()
